# Generate Report for Handback
#
# The localization status report is refreshed after a handback round-trip:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    for both locales (reflected on the Overview sheet and each locale sheet).
#  - Each locale sheet's "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns (J/K/L) are populated now that the
#    handback xliffs have come back for each source file.
#  - The "Latest Target File" cells become hyperlinks to the source .md file
#    on GitHub, matching the existing link style used in column A.
#  - A few columns are widened so the newly-populated long file names/links
#    are readable.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb98226b0bb0f910dda1d453126582c6fff566b0/e2e/"

$newStatus = "Handed back: in sync with en-US"

# -----------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# -----------------------------------------------------------------
# 2) zh-cn sheet: populate Latest Target File (J), Latest Handback File (K)
#    and Latest Handback DateTime (L) for both rows.
# -----------------------------------------------------------------
$zhRow2Name = "1ebf3363-8c0c-46ec-924f-82df1384e230.md"
$zhRow3Name = "7b9b633b-32b1-4527-ae6b-8e0e94fa10eb.md"

$wsZhCn.Range("J2").Value = $zhRow2Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), ($githubBase + $zhRow2Name), [Type]::Missing, [Type]::Missing, $zhRow2Name) | Out-Null
$wsZhCn.Range("K2").Value = "1ebf3363-8c0c-46ec-924f-82df1384e230.48041228292873ed592d0cf6afc48d10232dec1f.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-02-09 09:46:24"

$wsZhCn.Range("J3").Value = $zhRow3Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), ($githubBase + $zhRow3Name), [Type]::Missing, [Type]::Missing, $zhRow3Name) | Out-Null
$wsZhCn.Range("K3").Value = "7b9b633b-32b1-4527-ae6b-8e0e94fa10eb.4cfe42f02c9ab8e0a823fe9688e1e702d22c0ada.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-02-09 09:46:24"

# -----------------------------------------------------------------
# 3) de-de sheet: same columns, its own handback xliffs/datetime.
# -----------------------------------------------------------------
$deRow2Name = "1ebf3363-8c0c-46ec-924f-82df1384e230.md"
$deRow3Name = "7b9b633b-32b1-4527-ae6b-8e0e94fa10eb.md"

$wsDeDe.Range("J2").Value = $deRow2Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), ($githubBase + $deRow2Name), [Type]::Missing, [Type]::Missing, $deRow2Name) | Out-Null
$wsDeDe.Range("K2").Value = "1ebf3363-8c0c-46ec-924f-82df1384e230.48041228292873ed592d0cf6afc48d10232dec1f.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-02-09 09:46:53"

$wsDeDe.Range("J3").Value = $deRow3Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), ($githubBase + $deRow3Name), [Type]::Missing, [Type]::Missing, $deRow3Name) | Out-Null
$wsDeDe.Range("K3").Value = "7b9b633b-32b1-4527-ae6b-8e0e94fa10eb.4cfe42f02c9ab8e0a823fe9688e1e702d22c0ada.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-02-09 09:46:53"

# -----------------------------------------------------------------
# 4) Widen columns that now show long handback file names / links.
# -----------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15

$wsZhCn.Range("C1").ColumnWidth = 29.15
$wsZhCn.Range("J1").ColumnWidth = 39.15
$wsZhCn.Range("K1").ColumnWidth = 39.15

$wsDeDe.Range("C1").ColumnWidth = 29.15
$wsDeDe.Range("J1").ColumnWidth = 39.15
$wsDeDe.Range("K1").ColumnWidth = 39.15
